# Update on 12 Nov 2017
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Vendor Name"
$ws.Range("C1").Value = "Location From"
$ws.Range("D1").Value = "Location To"
$ws.Range("E1").Value = "KG"
$ws.Range("F1").Value = "Trips"

$ws.StandardWidth = 8
